$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 25 (shifts existing rows 25-78 down to 26-79)
$ws.Rows.Item(25).Insert()

# Populate the newly inserted row 25 with the new record
$ws.Range("A25").Value = 7
$ws.Range("B25").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C25").Value = "Ñuble"
$ws.Range("D25").Value = 44965
$ws.Range("E25").Value = 16
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103002
$ws.Range("J25").Value = "Ciruela"
$ws.Range("K25").Value = "Larry Ann"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 50
$ws.Range("N25").Value = 11000
$ws.Range("O25").Value = 11000
$ws.Range("P25").Value = 11000
$ws.Range("Q25").Value = "$/bandeja 18 kilos granel"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 611
$ws.Range("T25").Value = 18
